$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 78.818184
$ws.Range("I9").Value = 72
$ws.Range("K9").Value = 72
$ws.Range("M9").Value = 97
$ws.Range("H96").Value = 1961.2778
$ws.Range("J96").Value = 2817.182
$ws.Range("L96").Value = 8451.545999999998
$ws.Range("N96").Value = -11197.546
$ws.Range("H112").Value = 3256.3
$ws.Range("J112").Value = 3249.6
$ws.Range("L112").Value = 9748.799999999999
$ws.Range("N112").Value = -11964.8
$ws.Range("H137").Value = 4389.915
$ws.Range("I137").Value = 2125.975
$ws.Range("K137").Value = 6377.924999999999
$ws.Range("M137").Value = -3827.924999999999
$ws.Range("H138").Value = 4164.7803
$ws.Range("I138").Value = 3454.926
$ws.Range("K138").Value = 10364.778
$ws.Range("M138").Value = -5224.778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 244428.55
$ws.Range("I32").Value = 247951.19
$ws.Range("K32").Value = 247951.19
$ws.Range("M32").Value = -247664.19

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2603.0715
$ws.Range("I86").Value = 2524.25
$ws.Range("J86").Value = 2708.1667
$ws.Range("K86").Value = 2524.25
$ws.Range("L86").Value = 2708.1667
$ws.Range("M86").Value = -1401.25
$ws.Range("N86").Value = -4954.1667
$ws.Range("H89").Value = 2603.0715
$ws.Range("I89").Value = 2524.25
$ws.Range("J89").Value = 2708.1667
$ws.Range("K89").Value = 12621.25
$ws.Range("L89").Value = 13540.8335
$ws.Range("M89").Value = -7005.25
$ws.Range("N89").Value = -24772.8335
$ws.Range("H105").Value = 6166.615
$ws.Range("I105").Value = 2003.2222
$ws.Range("K105").Value = 2003.2222
$ws.Range("M105").Value = -256.2221999999999
$ws.Range("H139").Value = 81749
$ws.Range("J139").Value = 81749
$ws.Range("L139").Value = 81749
$ws.Range("N139").Value = -92029

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 83469.586
$ws.Range("J7").Value = 189.33333
$ws.Range("L7").Value = 189.33333
$ws.Range("N7").Value = -415.33333
$ws.Range("H31").Value = 3808.3
$ws.Range("I31").Value = 3811.8572
$ws.Range("J31").Value = 3800
$ws.Range("K31").Value = 3811.8572
$ws.Range("L31").Value = 3800
$ws.Range("M31").Value = -3516.8572
$ws.Range("N31").Value = -4390
$ws.Range("H34").Value = 3808.3
$ws.Range("I34").Value = 3811.8572
$ws.Range("J34").Value = 3800
$ws.Range("K34").Value = 3811.8572
$ws.Range("L34").Value = 3800
$ws.Range("M34").Value = -3609.8572
$ws.Range("N34").Value = -4204
$ws.Range("H94").Value = 5866.95
$ws.Range("J94").Value = 1049.5
$ws.Range("L94").Value = 1049.5
$ws.Range("N94").Value = -1951.5
$ws.Range("H99").Value = 27294.625
$ws.Range("I99").Value = 27294.625
$ws.Range("K99").Value = 27294.625
$ws.Range("M99").Value = -25796.625
$ws.Range("H107").Value = 893.0909
$ws.Range("I107").Value = 682.4
$ws.Range("K107").Value = 682.4
$ws.Range("M107").Value = 1237.6
$ws.Range("H126").Value = 27294.625
$ws.Range("I126").Value = 27294.625
$ws.Range("K126").Value = 81883.875
$ws.Range("M126").Value = -79413.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1916.2727
$ws.Range("I5").Value = 1545.5
$ws.Range("J5").Value = 2128.1428
$ws.Range("K5").Value = 4636.5
$ws.Range("L5").Value = 6384.428400000001
$ws.Range("M5").Value = -4524.5
$ws.Range("N5").Value = -6608.428400000001
$ws.Range("H92").Value = 893.625
$ws.Range("J92").Value = 992.7143
$ws.Range("L92").Value = 2978.1429
$ws.Range("N92").Value = -5474.1429
$ws.Range("H134").Value = 6550.2666
$ws.Range("I134").Value = 1006.75
$ws.Range("K134").Value = 3020.25
$ws.Range("M134").Value = 2049.75
$ws.Range("H135").Value = 1916.2727
$ws.Range("I135").Value = 1545.5
$ws.Range("J135").Value = 2128.1428
$ws.Range("K135").Value = 13909.5
$ws.Range("L135").Value = 19153.2852
$ws.Range("M135").Value = -11374.5
$ws.Range("N135").Value = -24223.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 3607231.5
$ws.Range("I2").Value = 5315829
$ws.Range("J2").Value = 192.66667
$ws.Range("K2").Value = 5315829
$ws.Range("L2").Value = 192.66667
$ws.Range("M2").Value = -5315716
$ws.Range("N2").Value = -418.66667
$ws.Range("H9").Value = 49
$ws.Range("I9").Value = 49
$ws.Range("K9").Value = 49
$ws.Range("M9").Value = 121
$ws.Range("H20").Value = 21453
$ws.Range("I20").Value = 18000
$ws.Range("K20").Value = 18000
$ws.Range("M20").Value = -17755
$ws.Range("H21").Value = 5680.3335
$ws.Range("I21").Value = 4816.9
$ws.Range("K21").Value = 4816.9
$ws.Range("M21").Value = -4643.9
$ws.Range("H24").Value = 412800
$ws.Range("I24").Value = 25000
$ws.Range("J24").Value = 671333.3
$ws.Range("K24").Value = 25000
$ws.Range("L24").Value = 671333.3
$ws.Range("M24").Value = -24827
$ws.Range("N24").Value = -671679.3
$ws.Range("H30").Value = 5680.3335
$ws.Range("I30").Value = 4816.9
$ws.Range("K30").Value = 4816.9
$ws.Range("M30").Value = -4711.9
$ws.Range("H54").Value = 18700
$ws.Range("J54").Value = 18700
$ws.Range("L54").Value = 18700
$ws.Range("N54").Value = -19480
$ws.Range("H132").Value = 17045.924
$ws.Range("I132").Value = 24206.938
$ws.Range("J132").Value = 5588.3
$ws.Range("K132").Value = 72620.814
$ws.Range("L132").Value = 16764.9
$ws.Range("M132").Value = -70090.814
$ws.Range("N132").Value = -21824.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 6000
$ws.Range("I3").Value = 6000
$ws.Range("K3").Value = 6000
$ws.Range("M3").Value = -5888
$ws.Range("H15").Value = 6000
$ws.Range("I15").Value = 6000
$ws.Range("K15").Value = 6000
$ws.Range("M15").Value = -5830
$ws.Range("H17").Value = 5000
$ws.Range("I17").Value = 5000
$ws.Range("K17").Value = 5000
$ws.Range("M17").Value = -4830
$ws.Range("H68").Value = 9395.962
$ws.Range("I68").Value = 7914.5
$ws.Range("K68").Value = 7914.5
$ws.Range("M68").Value = -7165.5
$ws.Range("H71").Value = 9395.962
$ws.Range("I71").Value = 7914.5
$ws.Range("K71").Value = 39572.5
$ws.Range("M71").Value = -35828.5
$ws.Range("H93").Value = 2369.2104
$ws.Range("I93").Value = 1002.6
$ws.Range("J93").Value = 7494
$ws.Range("K93").Value = 1002.6
$ws.Range("L93").Value = 7494
$ws.Range("M93").Value = 245.4
$ws.Range("N93").Value = -9990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 45000
$ws.Range("J21").Value = 50000
$ws.Range("L21").Value = 50000
$ws.Range("N21").Value = -50470
$ws.Range("H24").Value = 44722
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").Value = $null
$ws.Range("H35").Value = 45000
$ws.Range("J35").Value = 50000
$ws.Range("L35").Value = 50000
$ws.Range("N35").Value = -50580
